# Revert commit 6dc3ea4: restore the original Abstract/Authors text for rows
# 6, 10 and 11 (columns D/E) on Sheet1, undoing the extra-whitespace / dup
# shared-string edits that the reverted commit had introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = @"
id="Par1">The 2019 novel coronavirus disease (COVID-19) caused by severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) has spread globally, while the routes of transmission of this virus are still controversial.

 We enrolled 33 patients, without any ocular manifestation, with their ocular surface swabs collected for virus detection.

 RNA was detected strong positive in samples of both eyes from two patients.

 Therefore, SARS-CoV-2 may exist in the normal ocular surface of COVID-19 patients, suggesting that this virus might be spread through conjunctival contact.


"@

$ws.Range("E6").Value = @"
[Hua-Tao%Xie%NULL%1,   Shi-Yun%Jiang%NULL%1,   Kang-Kang%Xu%NULL%1,   Xin%Liu%NULL%1,   Bing%Xu%NULL%1,   Lin%Wang%lin_wang@hust.edu.cn%0,   Ming-Chang%Zhang%mingchangzhang@hotmail.com%1]
"@

$ws.Range("D10").Value = @"
Objectives
id="Par1">Since there are few reports on the ocular involvement of coronavirus disease 2019 (COVID-19) patients, this study aimed to assess the presence of severe acute respiratory syndrome Coronavirus-2 (SARS-CoV-2) in the tears of patients with COVID-19.
Methods
id="Par2">In this prospective case series, nasopharyngeal and tear sampling of 43 patients with severe COVID-19 were performed.

 The quantitative reverse transcription polymerase chain reaction (RT-PCR) was conducted to detect SARS-CoV-2. Ocular and systemic signs and symptoms were recorded from their medical history.


Results
id="Par3">The mean age of patients was 56 ± 13 years.

 The average disease time from initiation of symptoms was 3.27 days, range: 1–7 days.

 Forty-one patients (95.3%) had fever at the time of sampling.

 Only one patient had conjunctivitis.

 Thirty (69.8%) nasopharyngeal and three (7%) tear samples were positive for SARS-CoV-2. The result of tear sample was positive in the patient with conjunctivitis.

 All patients with positive tear RT-PCR results had positive nasopharyngeal RT-PCR results.


Conclusions
id="Par4">Ocular manifestation was rare in this series of severe COVID-19 patients, however, 7% of the patients had viral RNA in their conjunctival secretions.

 Therefore, possibility of ocular transmission should be considered even in the absence of ocular manifestations.



"@

$ws.Range("E10").Value = @"
[Saeed%Karimi%NULL%1,   Amir%Arabi%amir_arab_91@yahoo.com%1,   Toktam%Shahraki%NULL%2,   Toktam%Shahraki%NULL%0,   Sare%Safi%NULL%1]
"@

$ws.Range("D11").Value = @"
Background
id="Par1">A recent increase in children admitted with hypotensive shock and fever in the context of the COVID-19 outbreak requires an urgent characterization and assessment of the involvement of SARS-CoV-2 infection.

 This is a case series performed at 4 academic tertiary care centers in Paris of all the children admitted to the pediatric intensive care unit (PICU) with shock, fever and suspected SARS-CoV-2 infection between April 15th and April 27th, 2020.
Results
id="Par2">20 critically ill children admitted for shock had an acute myocarditis (left ventricular ejection fraction, 35% (25–55); troponin, 269 ng/mL (31–4607)), and arterial hypotension with mainly vasoplegic clinical presentation.

 The first symptoms before PICU admission were intense abdominal pain and fever for 6 days (1–10).

 All children had highly elevated C-reactive protein (&gt; 94 mg/L) and procalcitonin (&gt; 1.6 ng/mL) without microbial cause.

 At least one feature of Kawasaki disease was found in all children (fever, n = 20, skin rash, n = 10; conjunctivitis, n = 6; cheilitis, n = 5; adenitis, n = 2), but none had the typical form.

 SARS-CoV-2 PCR and serology were positive for 10 and 15 children, respectively.

 One child had both negative SARS-CoV-2 PCR and serology, but had a typical SARS-CoV-2 chest tomography scan.

 All children but one needed an inotropic/vasoactive drug support (epinephrine, n = 12; milrinone, n = 10; dobutamine, n = 6, norepinephrine, n = 4) and 8 were intubated.

 All children received intravenous immunoglobulin (2 g per kilogram) with adjuvant corticosteroids (n = 2), IL 1 receptor antagonist (n = 1) or a monoclonal antibody against IL-6 receptor (n = 1).

 All children survived and were afebrile with a full left ventricular function recovery at PICU discharge.


Conclusions
id="Par3">Acute myocarditis with intense systemic inflammation and atypical Kawasaki disease is an emerging severe pediatric disease following SARS-CoV-2 infection.

 Early recognition of this disease is needed and referral to an expert center is recommended.

 A delayed and inappropriate host immunological response is suspected.

 While underlying mechanisms remain unclear, further investigations are required to target an optimal treatment.



"@

$ws.Range("E11").Value = @"
[Marion%Grimaud%NULL%2,   Julie%Starck%NULL%2,   Michael%Levy%NULL%2,   Clémence%Marais%NULL%2,   Judith%Chareyre%NULL%2,   Diala%Khraiche%NULL%2,   Marianne%Leruez-Ville%NULL%2,   Pierre%Quartier%NULL%2,   Pierre Louis%Léger%NULL%2,   Guillaume%Geslain%NULL%2,   Nada%Semaan%NULL%2,   Florence%Moulin%NULL%2,   Matthieu%Bendavid%NULL%2,   Sandrine%Jean%NULL%2,   Géraldine%Poncelet%NULL%2,   Sylvain%Renolleau%NULL%2,   Mehdi%Oualha%mehdi.oualha@aphp.fr%3]
"@
